$d = $word.ActiveDocument

$d.Content.Find.Execute("64+32=96", $false, $true, $false, $false, $false, $true, 1, $false, "87-77=10", 2) | Out-Null
$d.Content.Find.Execute("6+4=10", $false, $true, $false, $false, $false, $true, 1, $false, "65-57=8", 2) | Out-Null
$d.Content.Find.Execute("77-29=48", $false, $true, $false, $false, $false, $true, 1, $false, "7+78=85", 2) | Out-Null
$d.Content.Find.Execute("58-53=5", $false, $true, $false, $false, $false, $true, 1, $false, "30+64=94", 2) | Out-Null
$d.Content.Find.Execute("3+76=79", $false, $true, $false, $false, $false, $true, 1, $false, "62-46=16", 2) | Out-Null
$d.Content.Find.Execute("36+2=38", $false, $true, $false, $false, $false, $true, 1, $false, "75-63=12", 2) | Out-Null
$d.Content.Find.Execute("51+7=58", $false, $true, $false, $false, $false, $true, 1, $false, "99-15=84", 2) | Out-Null
$d.Content.Find.Execute("63-33=30", $false, $true, $false, $false, $false, $true, 1, $false, "76+18=94", 2) | Out-Null
$d.Content.Find.Execute("2+23=25", $false, $true, $false, $false, $false, $true, 1, $false, "44-14=30", 2) | Out-Null
$d.Content.Find.Execute("64+26=90", $false, $true, $false, $false, $false, $true, 1, $false, "52+11=63", 2) | Out-Null
$d.Content.Find.Execute("44-13=31", $false, $true, $false, $false, $false, $true, 1, $false, "11+16=27", 2) | Out-Null
$d.Content.Find.Execute("22+36=58", $false, $true, $false, $false, $false, $true, 1, $false, "32-28=4", 2) | Out-Null
$d.Content.Find.Execute("95-77=18", $false, $true, $false, $false, $false, $true, 1, $false, "78-47=31", 2) | Out-Null
$d.Content.Find.Execute("34-25=9", $false, $true, $false, $false, $false, $true, 1, $false, "86-6=80", 2) | Out-Null
$d.Content.Find.Execute("42+17=59", $false, $true, $false, $false, $false, $true, 1, $false, "82-69=13", 2) | Out-Null
$d.Content.Find.Execute("77-20=57", $false, $true, $false, $false, $false, $true, 1, $false, "6-1=5", 2) | Out-Null
$d.Content.Find.Execute("16-7=9", $false, $true, $false, $false, $false, $true, 1, $false, "54+35=89", 2) | Out-Null
$d.Content.Find.Execute("23+61=84", $false, $true, $false, $false, $false, $true, 1, $false, "5+63=68", 2) | Out-Null
$d.Content.Find.Execute("33+29=62", $false, $true, $false, $false, $false, $true, 1, $false, "61-51=10", 2) | Out-Null
$d.Content.Find.Execute("70-57=13", $false, $true, $false, $false, $false, $true, 1, $false, "19+37=56", 2) | Out-Null
$d.Content.Find.Execute("81+12=93", $false, $true, $false, $false, $false, $true, 1, $false, "55-50=5", 2) | Out-Null
$d.Content.Find.Execute("35-30=5", $false, $true, $false, $false, $false, $true, 1, $false, "96-26=70", 2) | Out-Null
$d.Content.Find.Execute("31+33=64", $false, $true, $false, $false, $false, $true, 1, $false, "79-78=1", 2) | Out-Null
$d.Content.Find.Execute("13+9=22", $false, $true, $false, $false, $false, $true, 1, $false, "20+17=37", 2) | Out-Null
$d.Content.Find.Execute("43+15=58", $false, $true, $false, $false, $false, $true, 1, $false, "37+30=67", 2) | Out-Null
$d.Content.Find.Execute("4+88=92", $false, $true, $false, $false, $false, $true, 1, $false, "65-10=55", 2) | Out-Null
$d.Content.Find.Execute("8+37=45", $false, $true, $false, $false, $false, $true, 1, $false, "50+34=84", 2) | Out-Null
$d.Content.Find.Execute("36-28=8", $false, $true, $false, $false, $false, $true, 1, $false, "43+1=44", 2) | Out-Null
$d.Content.Find.Execute("42+36=78", $false, $true, $false, $false, $false, $true, 1, $false, "32-24=8", 2) | Out-Null
$d.Content.Find.Execute("39+55=94", $false, $true, $false, $false, $false, $true, 1, $false, "70-15=55", 2) | Out-Null
$d.Content.Find.Execute("67-65=2", $false, $true, $false, $false, $false, $true, 1, $false, "34+39=73", 2) | Out-Null
$d.Content.Find.Execute("23+0=23", $false, $true, $false, $false, $false, $true, 1, $false, "86+7=93", 2) | Out-Null
$d.Content.Find.Execute("60-21=39", $false, $true, $false, $false, $false, $true, 1, $false, "59-10=49", 2) | Out-Null
$d.Content.Find.Execute("14+70=84", $false, $true, $false, $false, $false, $true, 1, $false, "23+26=49", 2) | Out-Null
$d.Content.Find.Execute("77-57=20", $false, $true, $false, $false, $false, $true, 1, $false, "28+46=74", 2) | Out-Null
$d.Content.Find.Execute("12+71=83", $false, $true, $false, $false, $false, $true, 1, $false, "69-68=1", 2) | Out-Null
$d.Content.Find.Execute("26+15=41", $false, $true, $false, $false, $false, $true, 1, $false, "93-57=36", 2) | Out-Null
$d.Content.Find.Execute("47+43=90", $false, $true, $false, $false, $false, $true, 1, $false, "5+74=79", 2) | Out-Null
$d.Content.Find.Execute("13-12=1", $false, $true, $false, $false, $false, $true, 1, $false, "32+28=60", 2) | Out-Null
$d.Content.Find.Execute("23-20=3", $false, $true, $false, $false, $false, $true, 1, $false, "79-20=59", 2) | Out-Null
$d.Content.Find.Execute("97-83=14", $false, $true, $false, $false, $false, $true, 1, $false, "56+41=97", 2) | Out-Null
$d.Content.Find.Execute("72+9=81", $false, $true, $false, $false, $false, $true, 1, $false, "37+14=51", 2) | Out-Null
$d.Content.Find.Execute("51+9=60", $false, $true, $false, $false, $false, $true, 1, $false, "82-49=33", 2) | Out-Null
$d.Content.Find.Execute("14+40=54", $false, $true, $false, $false, $false, $true, 1, $false, "98-84=14", 2) | Out-Null
$d.Content.Find.Execute("35+17=52", $false, $true, $false, $false, $false, $true, 1, $false, "2+30=32", 2) | Out-Null
$d.Content.Find.Execute("25-21=4", $false, $true, $false, $false, $false, $true, 1, $false, "85-32=53", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $false, $true, $false, $false, $false, $true, 1, $false, "17-11=6", 2) | Out-Null
$d.Content.Find.Execute("74-33=41", $false, $true, $false, $false, $false, $true, 1, $false, "60-59=1", 2) | Out-Null
$d.Content.Find.Execute("90-36=54", $false, $true, $false, $false, $false, $true, 1, $false, "74+12=86", 2) | Out-Null
$d.Content.Find.Execute("71-32=39", $false, $true, $false, $false, $false, $true, 1, $false, "85-5=80", 2) | Out-Null
$d.Content.Find.Execute("46+51=97", $false, $true, $false, $false, $false, $true, 1, $false, "40-4=36", 2) | Out-Null
$d.Content.Find.Execute("2+14=16", $false, $true, $false, $false, $false, $true, 1, $false, "1+45=46", 2) | Out-Null
$d.Content.Find.Execute("99-63=36", $false, $true, $false, $false, $false, $true, 1, $false, "39-24=15", 2) | Out-Null
$d.Content.Find.Execute("40-34=6", $false, $true, $false, $false, $false, $true, 1, $false, "17+8=25", 2) | Out-Null
$d.Content.Find.Execute("38-17=21", $false, $true, $false, $false, $false, $true, 1, $false, "13+33=46", 2) | Out-Null
$d.Content.Find.Execute("34+4=38", $false, $true, $false, $false, $false, $true, 1, $false, "51+38=89", 2) | Out-Null
$d.Content.Find.Execute("6+85=91", $false, $true, $false, $false, $false, $true, 1, $false, "37+31=68", 2) | Out-Null
$d.Content.Find.Execute("77+16=93", $false, $true, $false, $false, $false, $true, 1, $false, "36+8=44", 2) | Out-Null
$d.Content.Find.Execute("19+77=96", $false, $true, $false, $false, $false, $true, 1, $false, "36+60=96", 2) | Out-Null
$d.Content.Find.Execute("38+41=79", $false, $true, $false, $false, $false, $true, 1, $false, "2+55=57", 2) | Out-Null
$d.Content.Find.Execute("65+32=97", $false, $true, $false, $false, $false, $true, 1, $false, "83+16=99", 2) | Out-Null
$d.Content.Find.Execute("40+10=50", $false, $true, $false, $false, $false, $true, 1, $false, "26+0=26", 2) | Out-Null
$d.Content.Find.Execute("8+85=93", $false, $true, $false, $false, $false, $true, 1, $false, "66+0=66", 2) | Out-Null
$d.Content.Find.Execute("81+17=98", $false, $true, $false, $false, $false, $true, 1, $false, "96-46=50", 2) | Out-Null
$d.Content.Find.Execute("58-5=53", $false, $true, $false, $false, $false, $true, 1, $false, "99-39=60", 2) | Out-Null
$d.Content.Find.Execute("77-24=53", $false, $true, $false, $false, $false, $true, 1, $false, "61-21=40", 2) | Out-Null
$d.Content.Find.Execute("3+54=57", $false, $true, $false, $false, $false, $true, 1, $false, "53-21=32", 2) | Out-Null
$d.Content.Find.Execute("20+78=98", $false, $true, $false, $false, $false, $true, 1, $false, "94-31=63", 2) | Out-Null
$d.Content.Find.Execute("43+42=85", $false, $true, $false, $false, $false, $true, 1, $false, "28-14=14", 2) | Out-Null
$d.Content.Find.Execute("22+25=47", $false, $true, $false, $false, $false, $true, 1, $false, "45-8=37", 2) | Out-Null
$d.Content.Find.Execute("39+50=89", $false, $true, $false, $false, $false, $true, 1, $false, "79-14=65", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $false, $true, $false, $false, $false, $true, 1, $false, "41-21=20", 2) | Out-Null
$d.Content.Find.Execute("87-44=43", $false, $true, $false, $false, $false, $true, 1, $false, "55+15=70", 2) | Out-Null
$d.Content.Find.Execute("35-29=6", $false, $true, $false, $false, $false, $true, 1, $false, "2+15=17", 2) | Out-Null
$d.Content.Find.Execute("15+70=85", $false, $true, $false, $false, $false, $true, 1, $false, "63-7=56", 2) | Out-Null
$d.Content.Find.Execute("80+0=80", $false, $true, $false, $false, $false, $true, 1, $false, "59+0=59", 2) | Out-Null
$d.Content.Find.Execute("1+53=54", $false, $true, $false, $false, $false, $true, 1, $false, "97-10=87", 2) | Out-Null
$d.Content.Find.Execute("8-8=0", $false, $true, $false, $false, $false, $true, 1, $false, "53+2=55", 2) | Out-Null
$d.Content.Find.Execute("52+7=59", $false, $true, $false, $false, $false, $true, 1, $false, "70-38=32", 2) | Out-Null
$d.Content.Find.Execute("72-47=25", $false, $true, $false, $false, $false, $true, 1, $false, "97-23=74", 2) | Out-Null
$d.Content.Find.Execute("68+27=95", $false, $true, $false, $false, $false, $true, 1, $false, "12+79=91", 2) | Out-Null
$d.Content.Find.Execute("2+64=66", $false, $true, $false, $false, $false, $true, 1, $false, "17+23=40", 2) | Out-Null
$d.Content.Find.Execute("98-5=93", $false, $true, $false, $false, $false, $true, 1, $false, "9+62=71", 2) | Out-Null
$d.Content.Find.Execute("76-68=8", $false, $true, $false, $false, $false, $true, 1, $false, "95+1=96", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $false, $true, $false, $false, $false, $true, 1, $false, "15+2=17", 2) | Out-Null
$d.Content.Find.Execute("23-2=21", $false, $true, $false, $false, $false, $true, 1, $false, "24-7=17", 2) | Out-Null
$d.Content.Find.Execute("10+41=51", $false, $true, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("57-29=28", $false, $true, $false, $false, $false, $true, 1, $false, "58+1=59", 2) | Out-Null
$d.Content.Find.Execute("28+57=85", $false, $true, $false, $false, $false, $true, 1, $false, "10+44=54", 2) | Out-Null
$d.Content.Find.Execute("31+27=58", $false, $true, $false, $false, $false, $true, 1, $false, "4+72=76", 2) | Out-Null
$d.Content.Find.Execute("73-63=10", $false, $true, $false, $false, $false, $true, 1, $false, "90-5=85", 2) | Out-Null
$d.Content.Find.Execute("25+60=85", $false, $true, $false, $false, $false, $true, 1, $false, "79+17=96", 2) | Out-Null
$d.Content.Find.Execute("81-48=33", $false, $true, $false, $false, $false, $true, 1, $false, "45+51=96", 2) | Out-Null
$d.Content.Find.Execute("69+5=74", $false, $true, $false, $false, $false, $true, 1, $false, "20+75=95", 2) | Out-Null
$d.Content.Find.Execute("60-8=52", $false, $true, $false, $false, $false, $true, 1, $false, "97-42=55", 2) | Out-Null
$d.Content.Find.Execute("97-34=63", $false, $true, $false, $false, $false, $true, 1, $false, "38-17=21", 2) | Out-Null
$d.Content.Find.Execute("8-4=4", $false, $true, $false, $false, $false, $true, 1, $false, "71-23=48", 2) | Out-Null
$d.Content.Find.Execute("51+22=73", $false, $true, $false, $false, $false, $true, 1, $false, "66+21=87", 2) | Out-Null
$d.Content.Find.Execute("99-82=17", $false, $true, $false, $false, $false, $true, 1, $false, "7+6=13", 2) | Out-Null
$d.Content.Find.Execute("72-41=31", $false, $true, $false, $false, $false, $true, 1, $false, "62+33=95", 2) | Out-Null
